$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (empty) column at N ---
# This shifts the previous N/O/P ("Late" / "heading" / "Outstanding") data
# one column to the right (-> O/P/Q), matching the new Variable-Instalments
# layout used by the RBI loan schedule.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()

# Give the freshly inserted column a sensible custom width (close to the
# neighbouring "In Advance" column) instead of the default width.
$ws.Columns("N").ColumnWidth = 9.9

# --- Make "Repayment schedule" the active/selected sheet & cell ---
# (was "NewLoanInput" before; now the 3rd tab, index 2, with J17 selected)
$ws.Activate()
$ws.Range("J17").Select() | Out-Null
